$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.789.94'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.85%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.339.67'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.35%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '546.36'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +6.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.91'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +4.74%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.562'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +7.02%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.38%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.152'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.357'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +8.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.755.04'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.63'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.85%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '57.758.69'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.83%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.43%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.345.98'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.63%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.61'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.96%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '334.74'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.23'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.67'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.58'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.98'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.167'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.28%  '
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.51'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.40'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +8.82%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +6.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '170.44'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0732'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +4.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.13'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.53%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.03'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +18.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.45'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.38%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.16'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +8.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.25'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.36%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +5.27%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '146.24'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.375'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.62'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '286.08'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0935'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0504'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '19.03'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +7.66%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.560'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.34%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.41%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.43'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +3.89%  '
